$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation row was inserted in front of the existing
# "2026/12/29" block, right after the last "2026/02/01" row (old row 752).
# Insert a whole new row at position 753, which pushes the old rows
# 753..794 down to 754..795 (values/format untouched) and extends the
# used range to A1:D795.
$ws.Rows.Item(753).Insert()

# Populate the newly inserted row 753 with the new data point.
# Prefix the date-like string with an apostrophe so Excel stores it as
# literal text (matching the rest of column A) instead of auto-converting
# it to a date serial number.
$ws.Range("A753").Value = "'2026/02/01"
$ws.Range("B753").Value = "日"
$ws.Range("C753").Value = 19
$ws.Range("D753").Value = 201
